# Generate Report for Handback
# - Marks the two e2e files as "Handed back: in sync with en-US" on the
#   Overview sheet as well as on each language sheet's Status column.
# - Fills in the "Latest Target File" / "Latest Handback File" /
#   "Latest Handback DateTime" columns for both rows on the zh-cn and
#   de-de sheets, turning the target-file cell into a hyperlink (same
#   target as the source-file-name hyperlink in column A).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: update the per-language status cells ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value2 = $statusText
$overview.Range("F2").Value2 = $statusText
$overview.Range("E3").Value2 = $statusText
$overview.Range("F3").Value2 = $statusText

# Overview status/target columns got wider to fit the longer text
# (~30 "characters" wide).
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# Source-file hyperlink targets (same for every language sheet -- they
# point at the source .md file in the repo).
$urlFile1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/5390735e348e7257c70a2ab2f3f48a211db4528d/e2e/2f65c9f5-f37b-45d5-973d-36e67bd6949b.md"
$urlFile2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/5390735e348e7257c70a2ab2f3f48a211db4528d/e2e/afb992b6-ec70-4159-8e56-3dcdcdf8286c.md"
$nameFile1 = "2f65c9f5-f37b-45d5-973d-36e67bd6949b.md"
$nameFile2 = "afb992b6-ec70-4159-8e56-3dcdcdf8286c.md"

function Update-LanguageSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C) on the language sheet.
    $ws.Range("C2").Value2 = $statusText
    $ws.Range("C3").Value2 = $statusText

    # Widen columns C (Status), I (Latest Target File) and J (Latest
    # Handback File) to fit the new, longer content.
    $ws.Columns.Item(3).ColumnWidth = 29.1666666666667
    $ws.Columns.Item(9).ColumnWidth = 39.1666666666667
    $ws.Columns.Item(10).ColumnWidth = 39.1666666666667

    # Rebuild the sheet's hyperlinks in order: A2, I2, A3, I3 -- this
    # keeps the existing A2/A3 links in place while inserting the new
    # I2/I3 ones right after their row's A-column link, matching how the
    # handback report generator lays relationships out.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $urlFile1, "", "", $nameFile1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlFile1, "", "", $nameFile1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlFile2, "", "", $nameFile2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlFile2, "", "", $nameFile2)

    # Latest Handback File (J) + Latest Handback DateTime (K): the
    # handback file reuses the same name as the Latest Handoff File (G)
    # for each row, and the datetime is the report-generation timestamp.
    $ws.Range("J2").Value2 = $ws.Range("G2").Value2
    $ws.Range("J3").Value2 = $ws.Range("G3").Value2

    $ws.Range("K2").Value2 = $HandbackDateTime
    $ws.Range("K3").Value2 = $HandbackDateTime
}

Update-LanguageSheet "zh-cn" "2016-08-13 02:58:14"
Update-LanguageSheet "de-de" "2016-08-13 02:58:24"
